$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string text "Surface Sites (lit) (mol/g)" -> "Surface Sites (lit) (mol)"
$cell = $ws.Cells.Find("Surface Sites (lit) (mol/g)")
if ($cell -ne $null) {
    $cell.Value = "Surface Sites (lit) (mol)"
}

# Convert the per-row formulas in columns E, G and H into shared formulas
# spanning E2:E7, G2:G7 and H2:H7 respectively (matching what Excel does
# when the same formula text is applied across a contiguous range).
$ws.Range("E2:E7").Formula = "=C2/D2"
$ws.Range("G2:G7").Formula = "=E2/F2"
$ws.Range("H2:H7").Formula = "=30/F2"

# Restore the current selection to B7 (previously B8)
$ws.Range("B7").Select()
